# Append the 4/4/2022 follow-up note to cell B4 on Sheet1.
# (Freeway econ prosperity report: running DeleteFeatures on the feature
#  layer of the line template instead of the feature class.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$existing = $ws.Range("B4").Value()

$newLine = "4/4/2022 - for one report (fwy econ prosperity), running DeleteFeatures on feature layer instead of feature class seems to fix it, but on other reports the issue doesn't seem to be cropping up."

$ws.Range("B4").Value = $existing + "`r`n" + $newLine
